# repeat Q factor data run for sg_rr_100_025 2023-12-11 14-23-14
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 93 (new, repeat run with prominence 1E-3) ----
$ws.Range("A93").Value = "sg_rr_100_025 2023-12-11 14-23-14.csv"
$ws.Range("B93").Value = 0.01
$ws.Range("C93").Value = 1000
$ws.Range("D93").Value = 5001
$ws.Range("E93").Value = 1530
$ws.Range("F93").Value = 1570
$ws.Range("G93").Value = 0.001
$ws.Range("H93").Value = "(approx_fsr/2)/wavelength step size"
$ws.Range("I93").Value = 1
$ws.Range("U93").Value = "seemed to find one peak in what looked like noise so increased prominence"

# ---- Row 94 (new, repeat run with prominence 1.5E-3, full result set) ----
$ws.Range("A94").Value = "sg_rr_100_025 2023-12-11 14-23-14.csv"
$ws.Range("B94").Value = 0.01
$ws.Range("C94").Value = 1000
$ws.Range("D94").Value = 5001
$ws.Range("E94").Value = 1530
$ws.Range("F94").Value = 1570
$ws.Range("G94").Value = 0.0015
$ws.Range("H94").Value = "(approx_fsr/2)/wavelength step size"
$ws.Range("I94").Value = 1
$ws.Range("J94").Value = 0.98282051282051597
$ws.Range("K94").Value = 0.0055097596875867197
$ws.Range("L94").Value = "yes although possibly missed peak at end, but hard to tell if that's really a peak anyway."
$ws.Range("M94").Value = 0.15507613665588901
$ws.Range("N94").Value = 0.0047121321070659299
$ws.Range("O94").Value = 10269.381068504999
$ws.Range("P94").Value = 236.72287624766
$ws.Range("Q94").Value = 108301085.330331
$ws.Range("R94").Value = 7492583.1937577203
$ws.Range("S94").Value = 100
$ws.Range("T94").Value = 0.1

# ---- Update the view: scroll position and active cell selection ----
$excel.ActiveWindow.ScrollRow = 78
$ws.Range("R94").Select()
